# Add a default header and a default footer to the document's (single)
# section, each containing one empty paragraph styled "Header"/"Footer"
# respectively - matching:
#   <w:headerReference w:type="default" r:id="rId9"/>
#   <w:footerReference w:type="default" r:id="rId10"/>
# plus the new word/header1.xml and word/footer1.xml parts.

$d = $word.ActiveDocument
$sec = $d.Sections.First

$hdr = $sec.Headers.Item(1)      # wdHeaderFooterPrimary -> "default" reference
$hdr.Range.set_Style("Header")

$ftr = $sec.Footers.Item(1)      # wdHeaderFooterPrimary -> "default" reference
$ftr.Range.set_Style("Footer")
